# Re-run NATMI for L1cam-Cd9 with the additional FAPs cluster per Dr Hou's advice.
# This refreshes the full results table (rows 2-17, cols A-T) for the
# 4x4 grid of Sending cluster x Target cluster combinations: ECs, M2, sCs, FAPs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A:T), reflecting the
# re-run of the NATMI ligand-receptor analysis (L1cam-Cd9) across all
# 4x4 Sending-cluster x Target-cluster combinations (ECs, M2, sCs, FAPs).
$rowsData = @(
    ,@("ECs","L1cam","Cd9","ECs",3.0,1.0,11.88712433333333,35.661373,0.5967229292030898,0.5967229292030898,3.0,1.0,40.314886,120.944658,0.1963336494301312,0.1963336494301312,479.2280623661593,4313.052561295433,0.1171567903890804,0.1171567903890804)
    ,@("ECs","L1cam","Cd9","M2",3.0,1.0,11.88712433333333,35.661373,0.5967229292030898,0.5967229292030898,3.0,1.0,69.65329366666667,208.959881,0.3392118072814421,0.3392118072814421,827.9773620418458,7451.796258376612,0.2024154632612561,0.2024154632612561)
    ,@("ECs","L1cam","Cd9","sCs",3.0,1.0,11.88712433333333,35.661373,0.5967229292030898,0.5967229292030898,3.0,1.0,69.40355333333332,208.21066,0.3379955709003387,0.3379955709003388,825.0086676484642,7425.078008836178,0.2016897071253207,0.2016897071253208)
    ,@("ECs","L1cam","Cd9","FAPs",3.0,1.0,11.88712433333333,35.661373,0.5967229292030898,0.5967229292030898,3.0,1.0,25.96691433333334,77.900743,0.126458972388088,0.126458972388088,308.6719392333488,2778.047453100139,0.07546096842743251,0.07546096842743251)
    ,@("M2","L1cam","Cd9","ECs",2.0,0.6666666666666666,0.3785896666666667,1.135769,0.01900486009268527,0.01900486009268527,3.0,1.0,40.314886,120.944658,0.1963336494301312,0.1963336494301312,15.26279925244467,137.365193272002,0.00373129353890596,0.00373129353890596)
    ,@("M2","L1cam","Cd9","M2",2.0,0.6666666666666666,0.3785896666666667,1.135769,0.01900486009268527,0.01900486009268527,3.0,1.0,69.65329366666667,208.959881,0.3392118072814421,0.3392118072814421,26.37001723149878,237.330155083489,0.006446672939170727,0.006446672939170727)
    ,@("M2","L1cam","Cd9","sCs",2.0,0.6666666666666666,0.3785896666666667,1.135769,0.01900486009268527,0.01900486009268527,3.0,1.0,69.40355333333332,208.21066,0.3379955709003387,0.3379955709003388,26.27546812194888,236.47921309754,0.006423558536908224,0.006423558536908225)
    ,@("M2","L1cam","Cd9","FAPs",2.0,0.6666666666666666,0.3785896666666667,1.135769,0.01900486009268527,0.01900486009268527,3.0,1.0,25.96691433333334,77.900743,0.126458972388088,0.126458972388088,9.830805441818557,88.47724897636701,0.002403335077700362,0.002403335077700362)
    ,@("sCs","L1cam","Cd9","ECs",3.0,1.0,4.123197999999999,12.369594,0.2069808238940481,0.2069808238940482,3.0,1.0,40.314886,120.944658,0.1963336494301312,0.1963336494301312,166.226257325428,1496.036315928852,0.04063730051717376,0.04063730051717377)
    ,@("sCs","L1cam","Cd9","M2",3.0,1.0,4.123197999999999,12.369594,0.2069808238940481,0.2069808238940482,3.0,1.0,69.65329366666667,208.959881,0.3392118072814421,0.3392118072814421,287.1943211398126,2584.748890258314,0.07021033934570196,0.07021033934570196)
    ,@("sCs","L1cam","Cd9","sCs",3.0,1.0,4.123197999999999,12.369594,0.2069808238940481,0.2069808238940482,3.0,1.0,69.40355333333332,208.21066,0.3379955709003387,0.3379955709003388,286.1645922968933,2575.481330672039,0.06995860173749127,0.0699586017374913)
    ,@("sCs","L1cam","Cd9","FAPs",3.0,1.0,4.123197999999999,12.369594,0.2069808238940481,0.2069808238940482,3.0,1.0,25.96691433333334,77.900743,0.126458972388088,0.126458972388088,107.0667292453713,963.600563208342,0.02617458229368114,0.02617458229368114)
    ,@("FAPs","L1cam","Cd9","ECs",3.0,1.0,3.531764333333333,10.595293,0.1772913868101768,0.1772913868101768,3.0,1.0,40.314886,120.944658,0.1963336494301312,0.1963336494301312,142.3826764771993,1281.444088294794,0.03480826498497102,0.03480826498497102)
    ,@("FAPs","L1cam","Cd9","M2",3.0,1.0,3.531764333333333,10.595293,0.1772913868101768,0.1772913868101768,3.0,1.0,69.65329366666667,208.959881,0.3392118072814421,0.3392118072814421,245.9990182711259,2213.991164440133,0.06013933173531328,0.06013933173531328)
    ,@("FAPs","L1cam","Cd9","sCs",3.0,1.0,3.531764333333333,10.595293,0.1772913868101768,0.1772913868101768,3.0,1.0,69.40355333333332,208.21066,0.3379955709003387,0.3379955709003388,245.1169942692644,2206.052948423379,0.05992370350061847,0.05992370350061849)
    ,@("FAPs","L1cam","Cd9","FAPs",3.0,1.0,3.531764333333333,10.595293,0.1772913868101768,0.1772913868101768,3.0,1.0,25.96691433333334,77.900743,0.126458972388088,0.126458972388088,91.70902188918879,825.3811970026991,0.02242008658927397,0.02242008658927397)
)

$nRows = $rowsData.Count
$nCols = $rowsData[0].Count
$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i,$j] = $rowsData[$i][$j]
    }
}

$ws.Range("A2:T17").Value2 = $arr
